$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The user selected the whole row 4 (Denis Nappa / denisnappa@gmail.com) and
# cleared its contents, leaving B4's hyperlink styling but no text/value,
# and removing the now-orphaned mailto hyperlink that pointed at B4.
$ws.Rows(4).Select() | Out-Null
$ws.Range("A4:B4").ClearContents() | Out-Null

foreach ($hl in @($ws.Hyperlinks)) {
    if ($hl.Range.Address($false, $false) -eq "B4") {
        $hl.Delete() | Out-Null
    }
}

$ws.Rows(4).Select() | Out-Null
